$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain text (matching the
# original inlineStr cell type); force text format before assigning so Excel
# does not auto-convert them to numbers.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D13", "D14", "D16", "D18", "D19", "D21", "D23", "D24", "D25", "D26", "D30", "D31", "D32", "D33", "D36", "D37", "D39", "D40", "D42", "D43", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.543.53"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.807.37"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "228.43"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "0.602"
$ws.Range("E6").Value = "  +7.66%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "36.53"
$ws.Range("E8").Value = "  +5.01%  "
$ws.Range("D9").Value = "0.300"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").Value = "0.0966"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "2.067.51"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "11.45"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.654"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.812.49"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "4.50"
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("D17").Value = "34.502.40"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "70.01"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "246.36"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D23").Value = "4.21"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +6.24%  "
$ws.Range("D25").Value = "172.81"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  +7.27%  "
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "4.02"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "3.85"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "0.0530"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "1.399.46"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("D36").Value = "0.672"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "2.44"
$ws.Range("E37").Value = "  -6.26%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "0.968"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "82.55"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  +7.19%  "
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "0.0495"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.968.61"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "104.30"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -2.39%  "
